$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column cells we are about to touch into Text format so that
# numeric-looking strings (e.g. "519.81", "1.00") are preserved verbatim instead
# of being coerced into numbers (which would drop formatting like trailing zeros
# or thousands separators).
$priceCells = @("D2","D3","D5","D6","D8","D9","D12","D13","D14","D15","D16","D17","D20","D21","D22","D23","D24","D26","D27","D28","D29","D30","D31","D33","D34","D35","D36","D37","D38","D39","D41","D42","D43","D44","D48","D49","D50","D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '72.324.37'
$ws.Range("E2").Value = '  +4.50%  '
$ws.Range("D3").Value = '4.034.23'
$ws.Range("E3").Value = '  +3.73%  '
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").Value = '519.81'
$ws.Range("E5").Value = '  -1.31%  '
$ws.Range("D6").Value = '147.31'
$ws.Range("E6").Value = '  +3.62%  '
$ws.Range("E7").Value = '  +2.49%  '
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  +0.15%  '
$ws.Range("D9").Value = '0.735'
$ws.Range("E9").Value = '  +2.54%  '
$ws.Range("E10").Value = '  +2.65%  '
$ws.Range("E11").Value = '  +0.97%  '
$ws.Range("D12").Value = '47.44'
$ws.Range("E12").Value = '  +13.30%  '
$ws.Range("D13").Value = '10.90'
$ws.Range("E13").Value = '  +7.24%  '
$ws.Range("D14").Value = '4.690.54'
$ws.Range("E14").Value = '  +4.09%  '
$ws.Range("D15").Value = '4.031.84'
$ws.Range("E15").Value = '  -0.02%  '
$ws.Range("D16").Value = '21.20'
$ws.Range("E16").Value = '  +7.97%  '
$ws.Range("D17").Value = '14.15'
$ws.Range("E17").Value = '  +2.92%  '
$ws.Range("E18").Value = '  -1.68%  '
$ws.Range("E19").Value = '  -2.07%  '
$ws.Range("D20").Value = '72.328.78'
$ws.Range("E20").Value = '  +4.67%  '
$ws.Range("D21").Value = '436.95'
$ws.Range("E21").Value = '  +3.15%  '
$ws.Range("D22").Value = '98.56'
$ws.Range("E22").Value = '  +12.37%  '
$ws.Range("D23").Value = '3.55'
$ws.Range("E23").Value = '  +6.35%  '
$ws.Range("D24").Value = '14.68'
$ws.Range("E24").Value = '  +3.99%  '
$ws.Range("E25").Value = '  -0.57%  '
$ws.Range("D26").Value = '11.88'
$ws.Range("E26").Value = '  +2.29%  '
$ws.Range("D27").Value = '11.26'
$ws.Range("E27").Value = '  +6.99%  '
$ws.Range("D28").Value = '37.32'
$ws.Range("E28").Value = '  +3.68%  '
$ws.Range("D29").Value = '3.07'
$ws.Range("E29").Value = '  +9.26%  '
$ws.Range("D30").Value = '13.50'
$ws.Range("E30").Value = '  +3.12%  '
$ws.Range("D31").Value = '690.33'
$ws.Range("E31").Value = '  -0.64%  '
$ws.Range("E32").Value = '  +2.78%  '
$ws.Range("D33").Value = '6.97'
$ws.Range("E33").Value = '  +17.90%  '
$ws.Range("D34").Value = '68.13'
$ws.Range("E34").Value = '  +0.43%  '
$ws.Range("D35").Value = '0.0₃0892'
$ws.Range("E35").Value = '  +6.94%  '
$ws.Range("B36").Value = 'ThetaToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D36").Value = '3.69'
$ws.Range("E36").Value = '  +25.22%  '
$ws.Range("B37").Value = 'TheGraph'
$ws.Range("C37").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D37").Value = '0.438'
$ws.Range("E37").Value = '  -0.72%  '
$ws.Range("D38").Value = '40.86'
$ws.Range("E38").Value = '  +1.93%  '
$ws.Range("D39").Value = '0.154'
$ws.Range("E39").Value = '  +3.68%  '
$ws.Range("E40").Value = '  -0.11%  '
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  +0.02%  '
$ws.Range("D42").Value = '0.0489'
$ws.Range("E42").Value = '  +2.26%  '
$ws.Range("D43").Value = '3.15'
$ws.Range("E43").Value = '  +4.99%  '
$ws.Range("D44").Value = '2.78'
$ws.Range("E44").Value = '  +0.54%  '
$ws.Range("E45").Value = '  +5.58%  '
$ws.Range("E46").Value = '  +4.53%  '
$ws.Range("E47").Value = '  +2.97%  '
$ws.Range("D48").Value = '9.04'
$ws.Range("E48").Value = '  +8.40%  '
$ws.Range("D49").Value = '0.000274'
$ws.Range("E49").Value = '  +21.68%  '
$ws.Range("D50").Value = '3.31'
$ws.Range("E50").Value = '  +1.32%  '
$ws.Range("D51").Value = '0.0₆0340'
$ws.Range("E51").Value = '  +1.81%  '

# Reset the number format back to the sheet default (General / style "Normal")
# now that the literal text has been stored, so no stray per-cell style survives.
foreach ($addr in $priceCells) {
    $ws.Range($addr).Style = "Normal"
}
